# Weekly update: insert the newest week's record at the top of the data
# table (row 531), pushing all existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 531, shifting rows 531:564
# down to 532:565 (and extending the sheet dimension to T565).
$ws.Rows("531:531").Insert()

# Populate the newly inserted row 531 with the latest price-report entry.
$ws.Range("A531").Value = 4
$ws.Range("B531").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C531").Value = "Los Lagos"
$ws.Range("D531").Value = 45267
$ws.Range("E531").Value = 10
$ws.Range("F531").Value = "Fruta"
$ws.Range("G531").Value = 100102
$ws.Range("H531").Value = "Cítricos"
$ws.Range("I531").Value = 100102004
$ws.Range("J531").Value = "Mandarina"
$ws.Range("K531").Value = "Murcott"
$ws.Range("L531").Value = "Primera"
$ws.Range("M531").Value = 300
$ws.Range("N531").Value = 7000
$ws.Range("O531").Value = 7000
$ws.Range("P531").Value = 7000
$ws.Range("Q531").Value = "$/bandeja 10 kilos"
$ws.Range("R531").Value = "Región de O'Higgins"
$ws.Range("S531").Value = 700
$ws.Range("T531").Value = 10
